$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.730.44'
$ws.Range('E2').Value = '  +0.97%  '
$ws.Range('D3').Value = '1.646.47'
$ws.Range('E3').Value = '  +1.18%  '
$ws.Range('E4').Value = '  +0.20%  '
$ws.Range('D5').Value = '215.94'
$ws.Range('E5').Value = '  +1.35%  '
$ws.Range('D6').Value = '0.507'
$ws.Range('E6').Value = '  +1.52%  '
$ws.Range('E7').Value = '  +0.18%  '
$ws.Range('D8').Value = '0.252'
$ws.Range('E8').Value = '  +1.37%  '
$ws.Range('E9').Value = '  +0.45%  '
$ws.Range('D10').Value = '19.19'
$ws.Range('E10').Value = '  +1.99%  '
$ws.Range('E11').Value = '  -0.24%  '
$ws.Range('D12').Value = '1.875.56'
$ws.Range('E12').Value = '  +1.18%  '
$ws.Range('E13').Value = '  +1.32%  '
$ws.Range('D14').Value = '1.603.38'
$ws.Range('E14').Value = '  -1.68%  '
$ws.Range('D15').Value = '0.532'
$ws.Range('E15').Value = '  +1.90%  '
$ws.Range('D16').Value = '65.33'
$ws.Range('E16').Value = '  +0.85%  '
$ws.Range('D17').Value = '26.724.19'
$ws.Range('E17').Value = '  +0.80%  '
$ws.Range('D18').Value = '0.0₃0744'
$ws.Range('E18').Value = '  +0.52%  '
$ws.Range('D19').Value = '218.41'
$ws.Range('E19').Value = '  +1.66%  '
$ws.Range('E20').Value = '  +0.14%  '
$ws.Range('D21').Value = '4.37'
$ws.Range('E21').Value = '  +1.53%  '
$ws.Range('E22').Value = '  +0.46%  '
$ws.Range('B23').Value = 'Toncoin'
$ws.Range('C23').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D23').Value = '2.29'
$ws.Range('E23').Value = '  +14.35%  '
$ws.Range('B24').Value = 'Avalanche'
$ws.Range('C24').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D24').Value = '9.53'
$ws.Range('E24').Value = '  +2.60%  '
$ws.Range('D25').Value = '145.90'
$ws.Range('E25').Value = '  -1.90%  '
$ws.Range('E26').Value = '  +0.13%  '
$ws.Range('E27').Value = '  +0.50%  '
$ws.Range('D28').Value = '7.13'
$ws.Range('E28').Value = '  +4.36%  '
$ws.Range('D29').Value = '15.74'
$ws.Range('E29').Value = '  +1.35%  '
$ws.Range('E30').Value = '  +1.46%  '
$ws.Range('E31').Value = '  +1.51%  '
$ws.Range('E32').Value = '  +1.15%  '
$ws.Range('D33').Value = '3.02'
$ws.Range('E33').Value = '  +2.44%  '
$ws.Range('D34').Value = '1.280.68'
$ws.Range('E34').Value = '  +5.15%  '
$ws.Range('E35').Value = '  +3.72%  '
$ws.Range('E36').Value = '  +1.86%  '
$ws.Range('D37').Value = '0.0179'
$ws.Range('E37').Value = '  +2.91%  '
$ws.Range('D38').Value = '0.537'
$ws.Range('E38').Value = '  +6.08%  '
$ws.Range('D39').Value = '0.826'
$ws.Range('E39').Value = '  +3.99%  '
$ws.Range('E40').Value = '  +0.16%  '
$ws.Range('D41').Value = '0.815'
$ws.Range('E41').Value = '  +2.87%  '
$ws.Range('E42').Value = '  -0.37%  '
$ws.Range('D43').Value = '5.45'
$ws.Range('E43').Value = '  +1.72%  '
$ws.Range('D44').Value = '1.785.46'
$ws.Range('E44').Value = '  +1.24%  '
$ws.Range('D45').Value = '92.00'
$ws.Range('E45').Value = '  -0.90%  '
$ws.Range('D46').Value = '59.83'
$ws.Range('E46').Value = '  +9.12%  '
$ws.Range('B48').Value = 'Cronos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D48').Value = '0.0515'
$ws.Range('E48').Value = '  +1.22%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').Value = '7.74'
$ws.Range('E49').Value = '  +3.29%  '
$ws.Range('B50').Value = 'Algorand'
$ws.Range('C50').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D50').Value = '0.0968'
$ws.Range('E50').Value = '  +1.90%  '
$ws.Range('B51').Value = 'Mantle'
$ws.Range('C51').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D51').Value = '0.406'
$ws.Range('E51').Value = '  -0.02%  '
